$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Status" header column right after the existing "Date" column
$ws.Range("G1").Value = "Status"

# Find the last used row of the table (header in row 1, data below)
$lastRow = $ws.UsedRange.Rows.Count()

# For every data row: refresh the Date column (F) and set the new Status (G)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "Feb 29, 2004 (12:00:00 EST)"
    $ws.Cells.Item($r, 7).Value = "Active"
}
